$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data previously in rows 2, 3, 4 is cyclically rotated:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
# Capture the original values first (Value2 for numbers, Text for strings)
# before overwriting any cells.

$cols = @("D","I","J","K","L","M","N","P","Q")

function Get-RowData($row) {
    $data = @{}
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $data[$col] = @{ Value = $cell.Value2; Text = $cell.Text }
    }
    return $data
}

$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row4 = Get-RowData 4

function Set-RowData($row, $data) {
    # Numeric columns get the numeric Value2; text columns (I and N) get the Text
    $ws.Range("D$row").Value = $data["D"].Value
    $ws.Range("I$row").Value = $data["I"].Text
    $ws.Range("J$row").Value = $data["J"].Value
    $ws.Range("K$row").Value = $data["K"].Value
    $ws.Range("L$row").Value = $data["L"].Value
    $ws.Range("M$row").Value = $data["M"].Value
    $ws.Range("N$row").Value = $data["N"].Text
    $ws.Range("P$row").Value = $data["P"].Value
    $ws.Range("Q$row").Value = $data["Q"].Value
}

# Apply rotation: row2 <- row4, row3 <- row2(old), row4 <- row3(old)
Set-RowData 2 $row4
Set-RowData 3 $row2
Set-RowData 4 $row3
